# The paragraph contains the run sequence:
#   ...david</pl>></m> Mays celuy...
# made of three adjacent runs:
#   1) "</pl>"  - Courier New 18, color 0000ff
#   2) ">"      - Arial 22, color 000000
#   3) "</m>"   - Courier New 18, color 0000ff
# The edit removes the middle ">" run entirely, which leaves the two
# Courier-New/blue runs adjacent to each other; Word (and this host)
# coalesces adjoining runs that share identical formatting, so runs
# 1 and 3 merge into a single run whose text is "</pl></m>" while
# keeping the original run's properties/rsid attributes.

$d = $word.ActiveDocument

$anchor = $d.Content
$found = $anchor.Find.Execute("david</pl>></m>", $true, $false, $false, $false, `
                               $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target 'david</pl>></m>' text"
}

$base = $anchor.Start
# Offset of the lone ">" character (the Arial run) within the matched text.
$gt = $d.Range($base + 10, $base + 11)

if ($gt.Text -ne ">") {
    throw "Unexpected text at computed offset: [$($gt.Text)]"
}

$gt.Delete()
